$wb = $excel.ActiveWorkbook

# ALC row 8: On the Drip | Eye Drops
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 26.818182
$ws.Range("I8").Value = 26.818182
$ws.Range("K8").Value = 80.45454599999999
$ws.Range("M8").Value = 58.54545400000001

# ALC row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 202.57143
$ws.Range("I28").Value = 249.63637
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 249.63637
$ws.Range("L28").Value = 30
$ws.Range("M28").Value = 235.36363
$ws.Range("N28").Value = -1000

# ALC row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8127.467
$ws.Range("I86").Value = 1602
$ws.Range("J86").Value = 17915.666
$ws.Range("K86").Value = 1602
$ws.Range("L86").Value = 17915.666
$ws.Range("M86").Value = -479
$ws.Range("N86").Value = -20161.666

# ALC row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 8127.467
$ws.Range("I89").Value = 1602
$ws.Range("J89").Value = 17915.666
$ws.Range("K89").Value = 8010
$ws.Range("L89").Value = 89578.33
$ws.Range("M89").Value = -2394
$ws.Range("N89").Value = -100810.33

# ARM row 12: Strait Ain't the Gate | Bronze Scutum
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 750
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -327
$ws.Range("N12").Value = -1346

# ARM row 63: Rivets Run through It | Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2606248
$ws.Range("I63").Value = 2374.125
$ws.Range("J63").Value = 7813996
$ws.Range("K63").Value = 2374.125
$ws.Range("L63").Value = 7813996
$ws.Range("M63").Value = -1688.125
$ws.Range("N63").Value = -7815368

# ARM row 66: A Riveting Revival (L) | Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2606248
$ws.Range("I66").Value = 2374.125
$ws.Range("J66").Value = 7813996
$ws.Range("K66").Value = 11870.625
$ws.Range("L66").Value = 39069980
$ws.Range("M66").Value = -8438.625
$ws.Range("N66").Value = -39076844

# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2273.7273
$ws.Range("I74").Value = 2254.926
$ws.Range("J74").Value = 2358.3333
$ws.Range("K74").Value = 2254.926
$ws.Range("L74").Value = 2358.3333
$ws.Range("M74").Value = -1380.926
$ws.Range("N74").Value = -4106.3333

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2273.7273
$ws.Range("I77").Value = 2254.926
$ws.Range("J77").Value = 2358.3333
$ws.Range("K77").Value = 11274.63
$ws.Range("L77").Value = 11791.6665
$ws.Range("M77").Value = -6906.629999999999
$ws.Range("N77").Value = -20527.6665

# ARM row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1377.4166
$ws.Range("I102").Value = 1169.8889
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1169.8889
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 452.1111000000001
$ws.Range("N102").Value = -5244

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1672.4286
$ws.Range("I107").Value = 1105.5
$ws.Range("J107").Value = 1899.2
$ws.Range("K107").Value = 1105.5
$ws.Range("L107").Value = 1899.2
$ws.Range("M107").Value = 814.5
$ws.Range("N107").Value = -5739.2

# BSM row 111: Heavy Hitter | Deepgold Knuckles
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 30633.334
$ws.Range("J111").Value = 30633.334
$ws.Range("L111").Value = 30633.334
$ws.Range("N111").Value = -38813.334

# CRP row 14: Citizens' Canes | Ash Radical
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 4999.6665
$ws.Range("J14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("N14").Value = -3840

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3011.95
$ws.Range("I31").Value = 1511.8334
$ws.Range("J31").Value = 5262.125
$ws.Range("K31").Value = 1511.8334
$ws.Range("L31").Value = 5262.125
$ws.Range("M31").Value = -1216.8334
$ws.Range("N31").Value = -5852.125

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3011.95
$ws.Range("I34").Value = 1511.8334
$ws.Range("J34").Value = 5262.125
$ws.Range("K34").Value = 1511.8334
$ws.Range("L34").Value = 5262.125
$ws.Range("M34").Value = -1309.8334
$ws.Range("N34").Value = -5666.125

# CRP row 99: O Pine | Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 16670030
$ws.Range("I99").Value = 3625576.5
$ws.Range("K99").Value = 3625576.5
$ws.Range("M99").Value = -3624078.5

# CRP row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 989.5
$ws.Range("I122").Value = 994.4375
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 2983.3125
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -533.3125
$ws.Range("N122").Value = -7750

# CRP row 126: A Better Conductor | Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 16670030
$ws.Range("I126").Value = 3625576.5
$ws.Range("K126").Value = 10876729.5
$ws.Range("M126").Value = -10874259.5

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2182.8518
$ws.Range("I132").Value = 1113.5555
$ws.Range("J132").Value = 4321.4443
$ws.Range("K132").Value = 3340.6665
$ws.Range("L132").Value = 12964.3329
$ws.Range("M132").Value = -810.6664999999998
$ws.Range("N132").Value = -18024.3329

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1112.6072
$ws.Range("I134").Value = 919.8
$ws.Range("J134").Value = 1335.0769
$ws.Range("K134").Value = 2759.4
$ws.Range("L134").Value = 4005.2307
$ws.Range("M134").Value = -224.3999999999996
$ws.Range("N134").Value = -9075.2307

# CUL row 12: Butter Me Up | Kukuru Butter
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 122
$ws.Range("J12").Value = 129.2
$ws.Range("L12").Value = 387.6
$ws.Range("N12").Value = -733.5999999999999

# CUL row 70: Persona non Gratin | Dhalmel Gratin
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3991.4285

# CUL row 73: Recipe for Disaster (L) | Dhalmel Gratin
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3991.4285

# CUL row 131: The Mountain Steeped | Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 766.35
$ws.Range("J131").Value = 783.9355
$ws.Range("L131").Value = 2351.8065
$ws.Range("N131").Value = -12431.8065

# GSM row 11: A Ringing Success | Copper Ring
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4662153
$ws.Range("I11").Value = 6500000
$ws.Range("J11").Value = 1721597.6
$ws.Range("K11").Value = 6500000
$ws.Range("L11").Value = 1721597.6
$ws.Range("M11").Value = -6499861
$ws.Range("N11").Value = -1721875.6

# GSM row 80: Needs More Prayerbell | Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4058.6667
$ws.Range("I80").Value = 3825
$ws.Range("J80").Value = 4175.5
$ws.Range("K80").Value = 3825
$ws.Range("L80").Value = 4175.5
$ws.Range("M80").Value = -2827
$ws.Range("N80").Value = -6171.5

# GSM row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4058.6667
$ws.Range("I83").Value = 3825
$ws.Range("J83").Value = 4175.5
$ws.Range("K83").Value = 19125
$ws.Range("L83").Value = 20877.5
$ws.Range("M83").Value = -14133
$ws.Range("N83").Value = -30861.5

# LTW row 22: Skin off Their Backs | Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5317
$ws.Range("J22").Value = 5325.25
$ws.Range("L22").Value = 5325.25
$ws.Range("N22").Value = -5915.25

# LTW row 27: Fire and Hide | Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5317
$ws.Range("J27").Value = 5325.25
$ws.Range("L27").Value = 5325.25
$ws.Range("N27").Value = -5539.25

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2371.353
$ws.Range("I132").Value = 1801
$ws.Range("J132").Value = 3186.1428
$ws.Range("K132").Value = 5403
$ws.Range("L132").Value = 9558.428400000001
$ws.Range("M132").Value = -2873
$ws.Range("N132").Value = -14618.4284

# LTW row 136: Respect for Br'aax | Br'aax Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1598.4
$ws.Range("I136").Value = 1444.8889
$ws.Range("K136").Value = 4334.6667
$ws.Range("M136").Value = -1784.6667

# WVR row 122: Heavy Armoire | Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1613.2963
$ws.Range("I122").Value = 1600.5769
$ws.Range("J122").Value = 1944
$ws.Range("K122").Value = 4801.7307
$ws.Range("L122").Value = 5832
$ws.Range("M122").Value = -2351.7307
$ws.Range("N122").Value = -10732
